$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "Sender Email"

# Sender email addresses, cycling every 4 rows down the gift list
$emails = @("jenny@gmail.com", "bobg@laterlly.co.uk", "jen@fal.com", "ban@dom.com")

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $email = $emails[$i % 4]
    $cell = $ws.Cells.Item($row, 3)
    $ws.Hyperlinks.Add($cell, "mailto:" + $email, "", "", $email)
}

$ws.Range("C9").Select() | Out-Null
